$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header columns to append after the existing K1 column
$ws.Range("L1").Value = "credit_limits.credit_limit [Currency]"
$ws.Range("M1").Value = "credit_limits.bypass_credit_limit_check [Check]"
$ws.Range("N1").Value = "sales_team.allocated_percentage [Float]"
$ws.Range("O1").Value = "sales_team.incentives [Currency]"

# Copy the header formatting (bold, centered, bordered) from the last
# existing header cell (K1) onto the newly added header cells.
$ws.Range("K1").Copy()
$ws.Range("L1:O1").PasteSpecial(-4122)
